$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows got reshuffled (rows rotated amongst each other) - likely from a
# "Fit" sorting/grouping bug fix. Apply the resulting row contents directly.

function Set-Row($r, $a, $b, $c, $d, $e, $f) {
    $ws.Cells.Item($r, 1).Value = $a
    $ws.Cells.Item($r, 2).Value = $b
    $ws.Cells.Item($r, 3).Value = $c
    $ws.Cells.Item($r, 4).Value = $d
    $ws.Cells.Item($r, 5).Value = $e
    $ws.Cells.Item($r, 6).Value = $f
}

Set-Row 3 901 16 15 45 60 60
Set-Row 4 1001 18 30 75 60 72
Set-Row 5 501 9 52 30 75 45
Set-Row 6 701 3 90 45 97 15

Set-Row 8 902 1 0 0 0 0
Set-Row 9 401 9 48 67 75 45
Set-Row 10 301 6 45 30 60 45

Set-Row 12 801 3 67 65 52 45
Set-Row 13 201 9 30 15 45 30
Set-Row 14 1201 2 10 10 10 10

Set-Row 16 3 0 3 3 3 3
Set-Row 17 1 0 2 2 2 2
Set-Row 18 1101 0 15 30 30 0

Set-Row 20 802 0 4 5 4 0
Set-Row 21 502 0 4 0 0 0
